# Commit: "Sun, Jul 05, 2020 11:05:06 AM"
#
# The only user-visible edit in this commit is a re-styling of the single
# table in the deck (slide 6) via PowerPoint's Table Styles gallery: the
# table's style id changes from
#   {6AF7B82C-4F09-438F-A33A-76C769BCF249}   (the deck's one custom/author
#                                              table style, "Table_0")
# to the built-in gallery style
#   {3C9DC34A-D94A-4ADB-AD50-882EAA1789AD}
#
# Table styles can only be changed through Table.ApplyStyle(id) - assigning
# Table.Style directly is rejected by the host ("Table styles cannot be
# assigned through a property").

$p = $ppt.ActivePresentation

$targetStyleId = "{3C9DC34A-D94A-4ADB-AD50-882EAA1789AD}"

$updated = $false
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $table = $shape.Table
            $table.ApplyStyle($targetStyleId)
            Write-Output "Slide $si / Shape $shi ('$($shape.Name)'): table style -> $($table.Style)"
            $updated = $true
        }
    }
}

if (-not $updated) {
    Write-Output "WARNING: no table shape found - no style was changed."
}
